{"js": "// The diff appends three new paragraphs to the end of the document body,\n// each one an exact duplicate of the existing title paragraph:\n//   text: \"Documento de Evidencias - DemoBlaze\"\n//   run formatting: font size 18pt (w:sz = 36 half-points)\nconst body = context.document.body;\n\nconst text = \"Documento de Evidencias - DemoBlaze\";\n\nfor (let i = 0; i < 3; i++) {\n  const paragraph = body.insertParagraph(text, Word.InsertLocation.end);\n  paragraph.font.size = 18;\n}\n\nawait context.sync();\n", "ps1": "# The diff appends three new paragraphs to the end of the document body,\n# each one an exact duplicate of the existing title paragraph:\n#   text: \"Documento de Evidencias - DemoBlaze\"\n#   run formatting: font size 18pt (w:sz = 36 half-points)\n$d = $word.ActiveDocument\n$text = \"Documento de Evidencias - DemoBlaze\"\n\nfor ($i = 0; $i -lt 3; $i++) {\n    $endRange = $d.Content\n    $endRange.Collapse(0)              # wdCollapseEnd\n    # Inserting a paragraph mark + the title text right after the current\n    # last run picks up that run's 18pt formatting, so the new paragraphs\n    # come out identical (same rPr/sz) to the ones already in the document.\n    $endRange.InsertAfter(\"`r\" + $text)\n}\n"}
